$wb = $excel.ActiveWorkbook

# Update zh-cn sheet (rows 4-7): Priority "low" -> "ht"; Latest Handoff Datetime updated
$wsZh = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-12-16 10:23:19"
}

# Update de-de sheet (rows 4-7): Priority "low" -> "ht"; Latest Handoff Datetime updated
$wsDe = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $wsDe.Cells.Item($r, 5).Value = "ht"
    $wsDe.Cells.Item($r, 8).Value = "2016-12-16 10:23:35"
}

# The de-de "Latest Handoff Datetime" (col H) shares its text with the Overview
# sheet's "Latest HO Xliff Generate Date" (col G) for these same rows, so that
# column needs the same refreshed value to stay in sync.
$wsOverview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-12-16 10:23:35"
}
